# Append a new log row (row 96) to each of the four "LIFTER" sheets,
# mirroring the previous row's static fields but with a fresh timestamp
# in column A. This mirrors the way the upstream logger appends a new
# sample to the database on each run.

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = 45772.44789538194
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x5a"
        E = "0xd"
        F = 400
        G = 568631262647113970876416.0
        H = 346
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = 45772.30395413194
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x56"
        E = "0xe"
        F = 400
        G = 568631262647113970876416.0
        H = 342
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = 45772.44768332176
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x5a"
        E = "0x3"
        F = 400
        G = 568631262647113970876416.0
        H = 346
        I = 3
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = 45772.51341248843
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x56"
        E = "0x3"
        F = 400
        G = 985046333984776000000000.0
        H = 342
        I = 3
    }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $ws.Cells.Item(96, 1).Value = $row.A
    $ws.Cells.Item(96, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item(96, 2).Value = $row.B
    $ws.Cells.Item(96, 3).Value = $row.C
    $ws.Cells.Item(96, 4).Value = $row.D
    $ws.Cells.Item(96, 5).Value = $row.E
    $ws.Cells.Item(96, 6).Value = $row.F
    $ws.Cells.Item(96, 7).Value = $row.G
    $ws.Cells.Item(96, 8).Value = $row.H
    $ws.Cells.Item(96, 9).Value = $row.I
}
